$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.710.56'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '3.166.53'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '613.39'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').Value = '146.26'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.162.34'
$ws.Range('E8').Value = '  +0.74%  '
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('D11').Value = '5.46'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '0.0000259'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').Value = '35.72'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '3.688.14'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('D17').Value = '64.702.97'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').Value = '3.165.75'
$ws.Range('E18').Value = '  +3.26%  '
$ws.Range('D19').Value = '6.88'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('D20').Value = '479.19'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').Value = '14.63'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '0.721'
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('D24').Value = '13.73'
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').Value = '84.15'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '8.75'
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('D28').Value = '2.81'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('D29').Value = '7.15'
$ws.Range('E29').Value = '  +4.31%  '
$ws.Range('E30').Value = '  -2.40%  '
$ws.Range('D31').Value = '2.11'
$ws.Range('E31').Value = '  -5.71%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('E33').Value = '  -1.15%  '
$ws.Range('D34').Value = '26.58'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('E36').Value = '  +7.70%  '
$ws.Range('D37').Value = '6.01'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = '53.32'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').Value = '3.19'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('D40').Value = '461.41'
$ws.Range('E40').Value = '  +2.03%  '
$ws.Range('D41').Value = '0.0399'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').Value = '8.34'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('D44').Value = '2.858.69'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('E45').Value = '  +2.89%  '
$ws.Range('D46').Value = '0.268'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('E47').Value = '  +5.49%  '
$ws.Range('D48').Value = '26.63'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').Value = '35.97'
$ws.Range('E50').Value = '  +7.63%  '
$ws.Range('D51').Value = '0.114'
$ws.Range('E51').Value = '  -0.17%  '
